$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Date/Time columns are stored as literal text, not converted
# into Excel date/time serial numbers.
$ws.Range("F2:G4").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "1cf4acd0-85b6-4bfd-967d-e2f3ff3c1165"
$ws.Range("B2").Value = "Out"
$ws.Range("C2").Value = "Tonore"
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = "2024-09-23"
$ws.Range("G2").Value = "18:11:24"

# Row 3
$ws.Range("A3").Value = "047a14ad-5d3d-4955-b6d1-4295cef09daf"
$ws.Range("B3").Value = "Out"
$ws.Range("C3").Value = "Ram"
$ws.Range("D3").Value = 120
$ws.Range("E3").Value = 19292
$ws.Range("F3").Value = "2024-09-23"
$ws.Range("G3").Value = "18:11:38"

# Row 4
$ws.Range("A4").Value = "482f3d17-ca96-473d-a88f-1cb2edb49201"
$ws.Range("B4").Value = "Out"
$ws.Range("C4").Value = "Ram"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1000
$ws.Range("F4").Value = "2024-09-23"
$ws.Range("G4").Value = "19:54:57"
